$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.225.20'
$ws.Range("E2").Value = '  +1.34%  '
$ws.Range("D3").Value = '1.859.65'
$ws.Range("E3").Value = '  +0.97%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.020'
$ws.Range("E5").Value = '  +1.40%  '
$ws.Range("D6").Value = '312.33'
$ws.Range("E6").Value = '  +0.96%  '
$ws.Range("D7").Value = '0.4798'
$ws.Range("E7").Value = '  +2.04%  '
$ws.Range("D8").Value = '0.3722'
$ws.Range("E8").Value = '  +1.72%  '
$ws.Range("D9").Value = '0.07323'
$ws.Range("E9").Value = '  +2.50%  '
$ws.Range("D10").Value = '0.9367'
$ws.Range("E10").Value = '  +1.11%  '
$ws.Range("D11").Value = '20.33'
$ws.Range("E11").Value = '  +4.02%  '
$ws.Range("D12").Value = '0.07863'
$ws.Range("E12").Value = '  +2.21%  '
$ws.Range("D13").Value = '1.848.25'
$ws.Range("E13").Value = '  -3.00%  '
$ws.Range("E14").Value = '  +2.52%  '
$ws.Range("D15").Value = '6.535'
$ws.Range("E15").Value = '  +2.21%  '
$ws.Range("D16").Value = '90.17'
$ws.Range("E16").Value = '  +2.31%  '
$ws.Range("E17").Value = '  +1.33%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008730'
$ws.Range("E18").Value = '  +1.24%  '
$ws.Range("E19").Value = '  +1.34%  '
$ws.Range("D20").Value = '14.76'
$ws.Range("E20").Value = '  +2.16%  '
$ws.Range("D21").Value = '27.257.06'
$ws.Range("E21").Value = '  +1.28%  '
$ws.Range("D22").Value = '5.105'
$ws.Range("E22").Value = '  +1.90%  '
$ws.Range("E23").Value = '  +0.74%  '
$ws.Range("D24").Value = '1.953'
$ws.Range("D25").Value = '153.81'
$ws.Range("E25").Value = '  +1.29%  '
$ws.Range("D26").Value = '18.51'
$ws.Range("E26").Value = '  +1.53%  '
$ws.Range("D27").Value = '1.997'
$ws.Range("E27").Value = '  -0.52%  '
$ws.Range("D28").Value = '115.68'
$ws.Range("E28").Value = '  +1.41%  '
$ws.Range("D29").Value = '4.986'
$ws.Range("E29").Value = '  +2.19%  '
$ws.Range("D30").Value = '0.08885'
$ws.Range("E30").Value = '  +0.71%  '
$ws.Range("D31").Value = '3.347'
$ws.Range("E31").Value = '  +4.10%  '
$ws.Range("E32").Value = '  +0.57%  '
$ws.Range("D33").Value = '4.584'
$ws.Range("E33").Value = '  +2.64%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7410'
$ws.Range("E34").Value = '  -0.74%  '
$ws.Range("D35").Value = '2.686'
$ws.Range("E35").Value = '  -3.42%  '
$ws.Range("D36").Value = '1.125'
$ws.Range("E36").Value = '  +3.67%  '
$ws.Range("D37").Value = '0.02033'
$ws.Range("E37").Value = '  +5.00%  '
$ws.Range("D38").Value = '0.05265'
$ws.Range("E38").Value = '  +1.24%  '
$ws.Range("D39").Value = '0.5327'
$ws.Range("E39").Value = '  +2.43%  '
$ws.Range("D40").Value = '7.107'
$ws.Range("E40").Value = '  +2.13%  '
$ws.Range("E41").Value = '  +1.27%  '
$ws.Range("D42").Value = '8.335'
$ws.Range("E42").Value = '  +2.30%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '10.60'
$ws.Range("E43").Value = '  +1.27%  '
$ws.Range("D44").Value = '0.4787'
$ws.Range("E44").Value = '  +1.96%  '
$ws.Range("D45").Value = '1.021'
$ws.Range("E45").Value = '  +1.45%  '
$ws.Range("D46").Value = '102.69'
$ws.Range("E46").Value = '  +1.26%  '
$ws.Range("D47").Value = '1.635'
$ws.Range("E47").Value = '  +2.57%  '
$ws.Range("D48").Value = '66.45'
$ws.Range("E48").Value = '  +1.17%  '
$ws.Range("D49").Value = '0.06078'
$ws.Range("E49").Value = '  +0.77%  '
$ws.Range("D50").Value = '0.9003'
$ws.Range("E50").Value = '  +1.12%  '
$ws.Range("D51").Value = '36.69'
$ws.Range("E51").Value = '  +1.28%  '
